$wb = $excel.ActiveWorkbook

# --- "Expected Out" sheet: update reconciled totals ---
$wsExpectedOut = $wb.Worksheets.Item("Expected Out")
$wsExpectedOut.Range("B9").Value = 1350.56
$wsExpectedOut.Range("B11").Value = 430.42
# B1 holds =SUM(B2:B295); Excel recalculates its cached value automatically.

# --- "Budget Out" sheet: update amount + description on row 9 ---
$wsBudgetOut = $wb.Worksheets.Item("Budget Out")
$wsBudgetOut.Range("C9").Value = 92.62
$wsBudgetOut.Range("F9").Value = "Description007zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

# --- "TestRecord" sheet: update date + amount on row 10 ---
$wsTestRecord = $wb.Worksheets.Item("TestRecord")
$wsTestRecord.Range("A10").Value = 43265
$wsTestRecord.Range("B10").Value = 123.54
$wsTestRecord.Range("E10").Value = "some test textzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

$excel.CalculateFull()
